# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.
#
# All source columns in this workbook are plain text (stored as inline
# strings in the original file) even though some values look like dates
# ("2026-02-01"), times ("19:59:19") or percentages ("77.0%"). Excel's
# normal `.Value =` assignment auto-converts those look-alike strings into
# real date/number values (and gives the cell a new number-format style in
# the process). To keep the appended rows as literal text - matching the
# rest of the sheet - each cell is briefly forced to Text format, written,
# then has its formatting cleared again so no stray style survives.
function Set-TextCell {
    param($cell, [string]$val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

function Add-LogRows {
    param($ws, [int]$startRow, $rows)
    $r = $startRow
    foreach ($row in $rows) {
        Set-TextCell $ws.Cells.Item($r, 1) $row[0]
        Set-TextCell $ws.Cells.Item($r, 2) $row[1]
        Set-TextCell $ws.Cells.Item($r, 3) $row[2]
        Set-TextCell $ws.Cells.Item($r, 4) $row[3]
        Set-TextCell $ws.Cells.Item($r, 5) $row[4]
        Set-TextCell $ws.Cells.Item($r, 6) $row[5]
        $r = $r + 1
    }
}

$wb = $excel.ActiveWorkbook

# ---- PIR sheet: append rows 83-96 -----------------------------------
$pirRows = @(
    ,@("2026-02-01","19:59:19","19:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","19:59:20","19:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","19:59:23","19:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","19:59:28","19:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","19:59:34","19:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","19:59:38","19:00","Bathroom","Motion Detected","Active")
    ,@("2026-02-01","19:59:46","19:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","19:59:51","19:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","19:59:56","19:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","20:00:01","20:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","20:00:06","20:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","20:00:11","20:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","20:00:16","20:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-01","20:00:18","20:00","Bathroom","Motion Detected","Active")
)
$wsPir = $wb.Worksheets.Item("PIR")
Add-LogRows $wsPir 83 $pirRows

# ---- Humidity sheet: append rows 65-75 -------------------------------
$humidityRows = @(
    ,@("2026-02-01","19:59:18","19:00","Bathroom","77.0%","Active")
    ,@("2026-02-01","19:59:19","19:00","Bathroom","77.8%","Active")
    ,@("2026-02-01","19:59:22","19:00","Bathroom","77.1%","Active")
    ,@("2026-02-01","19:59:27","19:00","Bathroom","78.0%","Active")
    ,@("2026-02-01","19:59:32","19:00","Bathroom","76.5%","Active")
    ,@("2026-02-01","19:59:37","19:00","Bathroom","77.4%","Active")
    ,@("2026-02-01","19:59:42","19:00","Bathroom","76.5%","Active")
    ,@("2026-02-01","19:59:47","19:00","Bathroom","77.8%","Active")
    ,@("2026-02-01","19:59:57","19:00","Bathroom","77.4%","Active")
    ,@("2026-02-01","20:00:07","20:00","Bathroom","76.9%","Active")
    ,@("2026-02-01","20:00:17","20:00","Bathroom","77.8%","Active")
)
$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRows $wsHumidity 65 $humidityRows

# ---- Temperature sheet: append rows 65-75 ----------------------------
$temperatureRows = @(
    ,@("2026-02-01","19:59:18","19:00","Bathroom","25.1C","Active")
    ,@("2026-02-01","19:59:19","19:00","Bathroom","25.1C","Active")
    ,@("2026-02-01","19:59:22","19:00","Bathroom","25.1C","Active")
    ,@("2026-02-01","19:59:27","19:00","Bathroom","25.1C","Active")
    ,@("2026-02-01","19:59:33","19:00","Bathroom","25.1C","Active")
    ,@("2026-02-01","19:59:37","19:00","Bathroom","25.1C","Active")
    ,@("2026-02-01","19:59:42","19:00","Bathroom","25.1C","Active")
    ,@("2026-02-01","19:59:47","19:00","Bathroom","25.1C","Active")
    ,@("2026-02-01","19:59:57","19:00","Bathroom","25.1C","Active")
    ,@("2026-02-01","20:00:08","20:00","Bathroom","25.1C","Active")
    ,@("2026-02-01","20:00:18","20:00","Bathroom","25.1C","Active")
)
$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-LogRows $wsTemperature 65 $temperatureRows
